# Updated cryptos list - refresh Price (column D) and Volume(1h) (column E)
# values for each coin row. A leading apostrophe is used for price values
# that would otherwise be auto-parsed by Excel as numbers, so they remain
# plain text exactly as scraped (preserving formats like "1.0000", "0.9999").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.425.95"
$ws.Range("E2").Value = "  -0.65%  "
$ws.Range("D3").Value = "1.725.11"
$ws.Range("E3").Value = "  -0.27%  "
$ws.Range("D4").Value = "'0.9995"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'243.40"
$ws.Range("E5").Value = "  -1.02%  "
$ws.Range("D6").Value = "'1.0000"
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("D8").Value = "'0.2607"
$ws.Range("E8").Value = "  -2.38%  "
$ws.Range("D9").Value = "'0.06205"
$ws.Range("E9").Value = "  +0.37%  "
$ws.Range("D10").Value = "1.715.05"
$ws.Range("E10").Value = "  -0.90%  "
$ws.Range("D11").Value = "'0.06987"
$ws.Range("E11").Value = "  -1.62%  "
$ws.Range("D12").Value = "'15.48"
$ws.Range("E12").Value = "  -0.98%  "
$ws.Range("D13").Value = "'4.544"
$ws.Range("E13").Value = "  +0.08%  "
$ws.Range("D14").Value = "'0.5998"
$ws.Range("E14").Value = "  -2.16%  "
$ws.Range("D15").Value = "'77.45"
$ws.Range("E15").Value = "  +0.25%  "
$ws.Range("D16").Value = "'0.9999"
$ws.Range("E16").Value = "  -0.01%  "
$ws.Range("D17").Value = "26.419.19"
$ws.Range("E17").Value = "  -0.64%  "
$ws.Range("D18").Value = "'0.9997"
$ws.Range("E18").Value = "  -0.05%  "
$ws.Range("D19").Value = "'0.000007223"
$ws.Range("E19").Value = "  +3.82%  "
$ws.Range("D20").Value = "'11.35"
$ws.Range("E20").Value = "  -1.74%  "
$ws.Range("D21").Value = "1.943.20"
$ws.Range("E21").Value = "  -0.62%  "
$ws.Range("D22").Value = "'4.463"
$ws.Range("E22").Value = "  -1.31%  "
$ws.Range("D23").Value = "'8.592"
$ws.Range("E23").Value = "  -2.49%  "
$ws.Range("D24").Value = "'5.159"
$ws.Range("E24").Value = "  -1.75%  "
$ws.Range("D25").Value = "'137.58"
$ws.Range("E25").Value = "  +0.12%  "
$ws.Range("E26").Value = "  -0.86%  "
$ws.Range("D27").Value = "'1.397"
$ws.Range("E27").Value = "  -0.99%  "
$ws.Range("D28").Value = "'106.94"
$ws.Range("E28").Value = "  -1.13%  "
$ws.Range("D29").Value = "'1.723"
$ws.Range("E29").Value = "  -3.04%  "
$ws.Range("D30").Value = "'3.955"
$ws.Range("E30").Value = "  -0.55%  "
$ws.Range("D31").Value = "'0.08002"
$ws.Range("E31").Value = "  -0.18%  "
$ws.Range("E32").Value = "  +0.07%  "
$ws.Range("D33").Value = "'0.04504"
$ws.Range("E33").Value = "  -1.03%  "
$ws.Range("D34").Value = "'0.9991"
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("E35").Value = "  -0.53%  "
$ws.Range("D36").Value = "'1.003"
$ws.Range("E36").Value = "  -0.22%  "
$ws.Range("D37").Value = "'0.6260"
$ws.Range("E37").Value = "  -1.00%  "
$ws.Range("D38").Value = "'0.9426"
$ws.Range("E38").Value = "  +4.96%  "
$ws.Range("D39").Value = "'2.389"
$ws.Range("E39").Value = "  +0.13%  "
$ws.Range("D40").Value = "'1.950"
$ws.Range("E40").Value = "  -4.65%  "
$ws.Range("D41").Value = "'0.9996"
$ws.Range("E41").Value = "  -0.17%  "
$ws.Range("D42").Value = "'0.01483"
$ws.Range("E42").Value = "  -1.20%  "
$ws.Range("D43").Value = "'99.64"
$ws.Range("E43").Value = "  -3.05%  "
$ws.Range("D44").Value = "'5.317"
$ws.Range("E44").Value = "  -2.14%  "
$ws.Range("D45").Value = "'0.3858"
$ws.Range("E45").Value = "  -1.04%  "
$ws.Range("D46").Value = "'6.835"
$ws.Range("E46").Value = "  -4.51%  "
$ws.Range("D47").Value = "'0.1171"
$ws.Range("E47").Value = "  -1.30%  "
$ws.Range("D48").Value = "'0.05364"
$ws.Range("E48").Value = "  -0.47%  "
$ws.Range("D49").Value = "'7.747"
$ws.Range("E49").Value = "  -1.65%  "
$ws.Range("D50").Value = "'30.23"
$ws.Range("E50").Value = "  -1.39%  "
$ws.Range("D51").Value = "'1.236"
$ws.Range("E51").Value = "  -1.43%  "